# ItemData.xlsx update:
#  - add a new item row (A16:D16) describing a "ground container" test asset
#  - fix the header row (A1:D1) so it shares the same (already existing,
#    identical) centered style record used elsewhere instead of a duplicate
#  - narrow column A a bit
#  - move the active selection to the newly added cell, matching the
#    author's last-saved cursor position

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- header row: reuse the existing centered style instead of the duplicate ---
# A1:D1 already use font "HarmonyOS Sans SC" (fontId 1); only the alignment
# needs to be (re)applied so the engine resolves it to the same cellXfs
# record already shared by the rest of the sheet (xlCenter = -4108).
$ws.Range("A1:D1").HorizontalAlignment = -4108
$ws.Range("A1:D1").VerticalAlignment = -4108

# --- new row 16: a ground-container test item ---
$ws.Range("A16").Value = 110000
$ws.Range("B16").Value = "地面容器"
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = "FLASE"

# --- column A a little narrower ---
$ws.Columns.Item(1).ColumnWidth = 7.4

# --- move selection to the newly added cell ---
$ws.Range("D16").Select() | Out-Null
